# Modified and updated all sleeper statements with setpage timeouts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Global Item" value for the data row from "Iphone" to "Cell Phones"
$ws.Range("C2").Value = "Cell Phones"

# Update the active selection on the sheet from H16 to E5
$ws.Range("E5").Select()
